$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '72.504.48'
$ws.Range("E2").Value = '  +1.66%  '
$ws.Range("D3").Value = '2.680.31'
$ws.Range("E3").Value = '  +2.09%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '''602.86'
$ws.Range("E5").Value = '  -0.62%  '
$ws.Range("D6").Value = '''178.64'
$ws.Range("E6").Value = '  -1.46%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '''0.525'
$ws.Range("E8").Value = '  +0.20%  '
$ws.Range("D9").Value = '2.678.93'
$ws.Range("E9").Value = '  +2.08%  '
$ws.Range("D10").Value = '''0.172'
$ws.Range("E10").Value = '  +3.99%  '
$ws.Range("E11").Value = '  +2.34%  '
$ws.Range("D13").Value = '''5.05'
$ws.Range("E13").Value = '  +0.50%  '
$ws.Range("D14").Value = '3.166.97'
$ws.Range("E14").Value = '  +2.89%  '
$ws.Range("D15").Value = '''0.0000187'
$ws.Range("E15").Value = '  +1.59%  '
$ws.Range("D16").Value = '72.386.83'
$ws.Range("E16").Value = '  +1.61%  '
$ws.Range("D17").Value = '''26.38'
$ws.Range("E17").Value = '  -0.85%  '
$ws.Range("D18").Value = '2.673.71'
$ws.Range("E18").Value = '  +2.04%  '
$ws.Range("D19").Value = '''11.93'
$ws.Range("E19").Value = '  +3.71%  '
$ws.Range("D20").Value = '''8.03'
$ws.Range("E20").Value = '  +1.37%  '
$ws.Range("D21").Value = '''373.60'
$ws.Range("E21").Value = '  -2.51%  '
$ws.Range("D22").Value = '''4.18'
$ws.Range("E22").Value = '  +1.24%  '
$ws.Range("D23").Value = '''2.05'
$ws.Range("E23").Value = '  +8.72%  '
$ws.Range("D24").Value = '''72.47'
$ws.Range("E24").Value = '  +0.17%  '
$ws.Range("E25").Value = '  -0.01%  '
$ws.Range("D26").Value = '''4.35'
$ws.Range("E26").Value = '  -2.62%  '
$ws.Range("D27").Value = '''9.89'
$ws.Range("E27").Value = '  +2.40%  '
$ws.Range("D28").Value = '2.815.09'
$ws.Range("E28").Value = '  +2.14%  '
$ws.Range("E29").Value = '  +0.13%  '
$ws.Range("D30").Value = '0.0₃0946'
$ws.Range("E30").Value = '  -1.50%  '
$ws.Range("B31").Value = 'Bittensor'
$ws.Range("C31").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D31").Value = '''518.77'
$ws.Range("E31").Value = '  -4.88%  '
$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").Value = '''8.09'
$ws.Range("E32").Value = '  +0.38%  '
$ws.Range("D33").Value = '''1.31'
$ws.Range("E33").Value = '  -1.18%  '
$ws.Range("E34").Value = '  -0.52%  '
$ws.Range("D36").Value = '''164.62'
$ws.Range("E36").Value = '  -0.59%  '
$ws.Range("D37").Value = '''19.54'
$ws.Range("E37").Value = '  +1.78%  '
$ws.Range("D38").Value = '''19.14'
$ws.Range("E38").Value = '  +0.85%  '
$ws.Range("D39").Value = '''1.39'
$ws.Range("E39").Value = '  +0.21%  '
$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").Value = '''1.84'
$ws.Range("E40").Value = '  -2.32%  '
$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").Value = '''0.110'
$ws.Range("E41").Value = '  -7.52%  '
$ws.Range("E42").Value = '  -0.01%  '
$ws.Range("D43").Value = '''5.05'
$ws.Range("E43").Value = '  -0.53%  '
$ws.Range("D44").Value = '''2.60'
$ws.Range("E44").Value = '  -2.36%  '
$ws.Range("D45").Value = '''0.334'
$ws.Range("E45").Value = '  +1.06%  '
$ws.Range("D46").Value = '''39.23'
$ws.Range("E46").Value = '  -2.25%  '
$ws.Range("D47").Value = '''154.20'
$ws.Range("E47").Value = '  -0.16%  '
$ws.Range("D48").Value = '''3.74'
$ws.Range("E48").Value = '  +2.71%  '
$ws.Range("D49").Value = '''0.549'
$ws.Range("E49").Value = '  +3.23%  '
$ws.Range("D50").Value = '''1.73'
$ws.Range("E50").Value = '  +2.19%  '
$ws.Range("D51").Value = '''0.0768'
$ws.Range("E51").Value = '  +1.74%  '
